$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Update the mass of the main 3D printed part (Pill_Puncher / Main file) to reflect the new CAD file
$ws.Range("D14").Value = 20

# Update the active selection to reflect where the user last clicked
$ws.Range("D14").Select()
